# Update the "btc2x" row (row 8) fee figures to reflect the most recent
# framework (fw) data, per commit message "update to recent fw".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 723.3
$ws.Range("D8").Value = 434
$ws.Range("E8").Value = 81729.3
$ws.Range("F8").Value = 54385.3
